# "FInal Commit for Project." -- add the final summary worksheet ("FInal Sheet")
# in front of the existing "Sheet1", containing the consolidated
# mean/variance/std-dev comparison table across the six tuning strategies.

$wb = $excel.ActiveWorkbook

# --- enable iterative calculation tolerance tweak (Options > Formulas) ----
$excel.MaxChange = 0.0001

# --- insert the new worksheet before the existing Sheet1 ------------------
$sheet1 = $wb.Worksheets.Item(1)
$ws = $wb.Worksheets.Add($sheet1)
$ws.Name = "FInal Sheet"

# Row labels down column A first (Sr No., the 15 sample indices, then the
# mean / var / StDev summary rows) so the shared-string table picks up
# "var" and "StDev" right after the already-known "Sr No." / "mean".
$ws.Range("A1").Value = "Sr No."
for ($i = 1; $i -le 15; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $i
}
$ws.Range("A17").Value = "mean"
$ws.Range("A18").Value = "var"
$ws.Range("A19").Value = "StDev"

# Column headers for the six result sets.
$ws.Range("B1").Value = "Pure GS"
$ws.Range("C1").Value = "PD H Reg."
$ws.Range("D1").Value = "GS H Reg."
$ws.Range("E1").Value = "PD H Log."
$ws.Range("F1").Value = "GS H Log."
$ws.Range("G1").Value = "Pure PD"

# Data body, rows 2-16.
$data = @(
    @(34.7,35.5,36.1,36.3,35,33.8),
    @(34.7,34.6,34.7,35.9,35.6,34.5),
    @(33.5,36.4,35.4,36.6,37,34.6),
    @(33.6,35,35.7,35.4,35.4,36),
    @(34.5,36.9,36.9,35.4,37.7,32.9),
    @(34.6,36.3,35.3,35.4,36.6,33.7),
    @(34.9,35.3,35.9,35,37.6,32.3),
    @(34.4,37.2,34.4,35.3,35.8,32.9),
    @(34.9,36.7,34.7,36,36.2,32),
    @(34.6,36.1,36.3,36.3,36.6,33.2),
    @(34.7,37.5,36.1,35.2,35.8,31.5),
    @(35.3,35.9,34.2,35.3,35.3,33.7),
    @(37.1,35,36.4,36.2,37.5,33.5),
    @(36,35.9,34.3,35.6,33.6,32.6),
    @(33.9,35.8,36.2,36.1,37.3,34.6)
)
for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 2).Value = $row[$c]
    }
}

# Summary formulas, rows 17-19.
$cols = @("B", "C", "D", "E", "F", "G")
foreach ($col in $cols) {
    $ws.Range($col + "17").Formula = "=AVERAGE(" + $col + "2:" + $col + "16)"
    $ws.Range($col + "18").Formula = "=VAR.S(" + $col + "2:" + $col + "16)"
    $ws.Range($col + "19").Formula = "=STDEV.S(" + $col + "2:" + $col + "16)"
}

$ws.Range("C7").Select()

$wb.Calculate()
